# Auto-generated edit script: applies numeric corrections to the Leve profit
# calculation columns (H-N) across multiple crafting-job sheets, per the
# scheduled runner's refreshed market-price snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1301.6522
$ws.Range("I15").Value = 1301.6522
$ws.Range("K15").Value = 3904.9566
$ws.Range("M15").Value = -3735.9566
$ws.Range("H18").Value = 6914.4
$ws.Range("I18").Value = 1842.4286
$ws.Range("J18").Value = 18749
$ws.Range("K18").Value = 1842.4286
$ws.Range("L18").Value = 18749
$ws.Range("M18").Value = -1558.4286
$ws.Range("N18").Value = -19317
$ws.Range("H62").Value = 286519.34
$ws.Range("I62").Value = 8325.777
$ws.Range("J62").Value = 453435.47
$ws.Range("K62").Value = 8325.777
$ws.Range("L62").Value = 453435.47
$ws.Range("M62").Value = -7701.777
$ws.Range("N62").Value = -454683.47
$ws.Range("H65").Value = 286519.34
$ws.Range("I65").Value = 8325.777
$ws.Range("J65").Value = 453435.47
$ws.Range("K65").Value = 41628.885
$ws.Range("L65").Value = 2267177.35
$ws.Range("M65").Value = -38508.885
$ws.Range("N65").Value = -2273417.35
$ws.Range("H69").Value = 21750
$ws.Range("I69").Value = 27800
$ws.Range("J69").Value = 15700
$ws.Range("K69").Value = 83400
$ws.Range("L69").Value = 47100
$ws.Range("M69").Value = -82526
$ws.Range("N69").Value = -48848
$ws.Range("H72").Value = 21750
$ws.Range("I72").Value = 27800
$ws.Range("J72").Value = 15700
$ws.Range("K72").Value = 250200
$ws.Range("L72").Value = 141300
$ws.Range("M72").Value = -245832
$ws.Range("N72").Value = -150036
$ws.Range("H76").Value = 125003130
$ws.Range("J76").Value = 4999
$ws.Range("L76").Value = 4999
$ws.Range("N76").Value = -5629
$ws.Range("H79").Value = 125003130
$ws.Range("J79").Value = 4999
$ws.Range("L79").Value = 4999
$ws.Range("N79").Value = -7183
$ws.Range("H100").Value = 2209.1765
$ws.Range("I100").Value = 1795.9
$ws.Range("J100").Value = 2799.5715
$ws.Range("K100").Value = 1795.9
$ws.Range("L100").Value = 2799.5715
$ws.Range("M100").Value = -1254.9
$ws.Range("N100").Value = -3881.5715
$ws.Range("H106").Value = 508557.53
$ws.Range("I106").Value = 928255
$ws.Range("J106").Value = 4920.6
$ws.Range("K106").Value = 928255
$ws.Range("L106").Value = 4920.6
$ws.Range("M106").Value = -927624
$ws.Range("N106").Value = -6182.6
$ws.Range("H135").Value = 1214
$ws.Range("I135").Value = 1268.6364
$ws.Range("K135").Value = 11417.7276
$ws.Range("M135").Value = -8882.7276
$ws.Range("H138").Value = 2426.6038
$ws.Range("I138").Value = 1680.1538
$ws.Range("J138").Value = 2669.2
$ws.Range("K138").Value = 5040.4614
$ws.Range("L138").Value = 8007.599999999999
$ws.Range("M138").Value = 99.53859999999986
$ws.Range("N138").Value = -18287.6
$ws.Range("H141").Value = 1418.7142
$ws.Range("I141").Value = 1571.8334
$ws.Range("K141").Value = 4715.5002
$ws.Range("M141").Value = 464.4997999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1385.625
$ws.Range("I2").Value = 1244.6154
$ws.Range("K2").Value = 1244.6154
$ws.Range("M2").Value = -1131.6154
$ws.Range("H32").Value = 12796.487
$ws.Range("I32").Value = 7829.8335
$ws.Range("K32").Value = 7829.8335
$ws.Range("M32").Value = -7542.8335
$ws.Range("H45").Value = 3036
$ws.Range("I45").Value = 2343.7144
$ws.Range("J45").Value = 4247.5
$ws.Range("K45").Value = 2343.7144
$ws.Range("L45").Value = 4247.5
$ws.Range("M45").Value = -1966.7144
$ws.Range("N45").Value = -5001.5
$ws.Range("H116").Value = 1385.625
$ws.Range("I116").Value = 1244.6154
$ws.Range("K116").Value = 1244.6154
$ws.Range("M116").Value = 1049.3846
$ws.Range("H139").Value = 74047.336
$ws.Range("J139").Value = 74047.336
$ws.Range("L139").Value = 74047.336
$ws.Range("N139").Value = -84327.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1385.625
$ws.Range("I3").Value = 1244.6154
$ws.Range("K3").Value = 1244.6154
$ws.Range("M3").Value = -1130.6154
$ws.Range("H36").Value = 1600
$ws.Range("I36").Value = 1600
$ws.Range("K36").Value = 1600
$ws.Range("M36").Value = -1066
$ws.Range("H107").Value = 3620.55
$ws.Range("I107").Value = 2785.2307
$ws.Range("J107").Value = 5171.857
$ws.Range("K107").Value = 2785.2307
$ws.Range("L107").Value = 5171.857
$ws.Range("M107").Value = -865.2307000000001
$ws.Range("N107").Value = -9011.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1728.0588
$ws.Range("I16").Value = 1412.7142
$ws.Range("J16").Value = 3199.6667
$ws.Range("K16").Value = 1412.7142
$ws.Range("L16").Value = 3199.6667
$ws.Range("M16").Value = -1125.7142
$ws.Range("N16").Value = -3773.6667
$ws.Range("H31").Value = 3125.6667
$ws.Range("I31").Value = 2258.5
$ws.Range("K31").Value = 2258.5
$ws.Range("M31").Value = -1963.5
$ws.Range("H34").Value = 3125.6667
$ws.Range("I34").Value = 2258.5
$ws.Range("K34").Value = 2258.5
$ws.Range("M34").Value = -2056.5
$ws.Range("H62").Value = 5518.8667
$ws.Range("I62").Value = 4989.636
$ws.Range("J62").Value = 6974.25
$ws.Range("K62").Value = 4989.636
$ws.Range("L62").Value = 6974.25
$ws.Range("M62").Value = -4365.636
$ws.Range("N62").Value = -8222.25
$ws.Range("H65").Value = 5518.8667
$ws.Range("I65").Value = 4989.636
$ws.Range("J65").Value = 6974.25
$ws.Range("K65").Value = 24948.18
$ws.Range("L65").Value = 34871.25
$ws.Range("M65").Value = -21828.18
$ws.Range("N65").Value = -41111.25
$ws.Range("H94").Value = 829.931
$ws.Range("I94").Value = 688.44446
$ws.Range("J94").Value = 893.6
$ws.Range("K94").Value = 688.44446
$ws.Range("L94").Value = 893.6
$ws.Range("M94").Value = -237.44446
$ws.Range("N94").Value = -1795.6
$ws.Range("H113").Value = 1728.0588
$ws.Range("I113").Value = 1412.7142
$ws.Range("J113").Value = 3199.6667
$ws.Range("K113").Value = 1412.7142
$ws.Range("L113").Value = 3199.6667
$ws.Range("M113").Value = 757.2858000000001
$ws.Range("N113").Value = -7539.6667
$ws.Range("H132").Value = 2041.2963
$ws.Range("I132").Value = 1864.1052
$ws.Range("K132").Value = 5592.3156
$ws.Range("M132").Value = -3062.3156
$ws.Range("H134").Value = 36973.32
$ws.Range("I134").Value = 1063.2916
$ws.Range("K134").Value = 3189.8748
$ws.Range("M134").Value = -654.8748000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1511.4117
$ws.Range("I5").Value = 1147.3
$ws.Range("J5").Value = 2031.5714
$ws.Range("K5").Value = 3441.9
$ws.Range("L5").Value = 6094.7142
$ws.Range("M5").Value = -3329.9
$ws.Range("N5").Value = -6318.7142
$ws.Range("H59").Value = 1260
$ws.Range("J59").Value = 3000
$ws.Range("L59").Value = 9000
$ws.Range("N59").Value = -10080
$ws.Range("H60").Value = 877.9783
$ws.Range("I60").Value = 871.4
$ws.Range("K60").Value = 2614.2
$ws.Range("M60").Value = -2363.2
$ws.Range("H61").Value = 84.666664
$ws.Range("I61").Value = 103.42857
$ws.Range("J61").Value = 19
$ws.Range("K61").Value = 310.28571
$ws.Range("L61").Value = 57
$ws.Range("M61").Value = -95.28570999999999
$ws.Range("N61").Value = -487
$ws.Range("H118").Value = 3166.8
$ws.Range("I118").Value = 839
$ws.Range("J118").Value = 3748.75
$ws.Range("K118").Value = 2517
$ws.Range("L118").Value = 11246.25
$ws.Range("M118").Value = -1274
$ws.Range("N118").Value = -13732.25
$ws.Range("H133").Value = 5999
$ws.Range("I133").Value = 3998.5
$ws.Range("J133").Value = 10000
$ws.Range("K133").Value = 11995.5
$ws.Range("L133").Value = 30000
$ws.Range("M133").Value = -6935.5
$ws.Range("N133").Value = -40120
$ws.Range("H134").Value = 899
$ws.Range("I134").Value = 899
$ws.Range("K134").Value = 2697
$ws.Range("M134").Value = 2373
$ws.Range("H135").Value = 1511.4117
$ws.Range("I135").Value = 1147.3
$ws.Range("J135").Value = 2031.5714
$ws.Range("K135").Value = 10325.7
$ws.Range("L135").Value = 18284.1426
$ws.Range("M135").Value = -7790.699999999999
$ws.Range("N135").Value = -23354.1426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 291059.25
$ws.Range("I70").Value = 270462.34
$ws.Range("K70").Value = 270462.34
$ws.Range("M70").Value = -270192.34
$ws.Range("H73").Value = 291059.25
$ws.Range("I73").Value = 270462.34
$ws.Range("K73").Value = 270462.34
$ws.Range("M73").Value = -269526.34

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 957
$ws.Range("J22").Value = 1070.0769
$ws.Range("L22").Value = 1070.0769
$ws.Range("N22").Value = -1660.0769
$ws.Range("H27").Value = 957
$ws.Range("J27").Value = 1070.0769
$ws.Range("L27").Value = 1070.0769
$ws.Range("N27").Value = -1284.0769
$ws.Range("H46").Value = 1519.6
$ws.Range("I46").Value = 649.5
$ws.Range("K46").Value = 649.5
$ws.Range("M46").Value = -461.5
$ws.Range("H68").Value = 2878
$ws.Range("I68").Value = 2878
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2878
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -2129
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 2878
$ws.Range("I71").Value = 2878
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 14390
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -10646
$ws.Range("N71").ClearContents()
$ws.Range("H136").Value = 4220.5557
$ws.Range("I136").Value = 4396.4116
$ws.Range("J136").Value = 3921.6
$ws.Range("K136").Value = 13189.2348
$ws.Range("L136").Value = 11764.8
$ws.Range("M136").Value = -10639.2348
$ws.Range("N136").Value = -16864.8
